$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the separate "Author" byline paragraph ("Last updated by
#    Randy Peterson"). The remaining "Date" paragraph becomes the
#    sole byline line (its style was already "Date").
# ------------------------------------------------------------------
$authorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Style.NameLocal -eq "Author") {
        $authorIndex = $i
        break
    }
}
if ($authorIndex -gt 0) {
    $d.Paragraphs($authorIndex).Range.Delete()
}

# ------------------------------------------------------------------
# 2) Bump the last-updated date from "October 2, 2018" to
#    "October 3, 2018" inside the (now single) "Date" paragraph.
#    We rebuild the trailing "2, 2018" portion run-by-run (delete +
#    append) instead of doing an in-place text replace, so the
#    untouched runs ("Last", " ", "updated", " ", "on", " ",
#    "October", " ") keep their original identities and only the
#    day number actually changes.
# ------------------------------------------------------------------
$dateIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Style.NameLocal -eq "Date") {
        $dateIndex = $i
        break
    }
}

if ($dateIndex -gt 0) {
    $datePara = $d.Paragraphs($dateIndex)
    $find = $datePara.Range
    $find.Find.Execute("2, 2018", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $find.Delete()

    $pEnd = $datePara.Range.End
    $ins = $d.Range($pEnd - 1, $pEnd - 1)
    $ins.InsertAfter("3,")

    $pEnd = $datePara.Range.End
    $ins = $d.Range($pEnd - 1, $pEnd - 1)
    $ins.InsertAfter(" ")

    $pEnd = $datePara.Range.End
    $ins = $d.Range($pEnd - 1, $pEnd - 1)
    $ins.InsertAfter("2018")
}
